$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.252.52"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "1.590.40"
$ws.Range("E3").Value = "  +1.22%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.18"
$ws.Range("E5").Value = "  +1.73%  "
$ws.Range("E6").Value = "  +1.01%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  +0.65%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.38"
$ws.Range("E10").Value = "  -0.95%  "
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("D12").Value = "1.814.02"
$ws.Range("E12").Value = "  +1.22%  "
$ws.Range("D13").Value = "1.586.38"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.37"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").Value = "26.255.10"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.42"
$ws.Range("E19").Value = "  +2.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "213.67"
$ws.Range("E20").Value = "  +3.27%  "
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("E22").Value = "  +1.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.02"
$ws.Range("E23").Value = "  +1.96%  "
$ws.Range("E24").Value = "  -2.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.47"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("E27").Value = "  +1.46%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  -1.17%  "
$ws.Range("E31").Value = "  +1.48%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.94"
$ws.Range("E33").Value = "  -0.88%  "
$ws.Range("D34").Value = "1.338.61"
$ws.Range("E34").Value = "  +4.64%  "
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.592"
$ws.Range("E37").Value = "  -2.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0166"
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("E39").Value = "  +0.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.78"
$ws.Range("E40").Value = "  +3.86%  "
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -6.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.14"
$ws.Range("E43").Value = "  +0.58%  "
$ws.Range("E44").Value = "  +0.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.86"
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("D46").Value = "1.725.68"
$ws.Range("E46").Value = "  +1.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "85.62"
$ws.Range("E48").Value = "  -2.48%  "
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0976"
$ws.Range("E50").Value = "  -2.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("E51").Value = "  -0.31%  "
